{"js": "// The Matlab \"pyenv(...)\" code sample paragraph ended with a stray\n// trailing double-quote character (a leftover from copy/pasting a\n// console transcript): ...'Scripts', 'python.EXE'))\"\n// The revision removes that stray closing quote so the line reads\n// ...'Scripts', 'python.EXE')) -- matching the other, correctly\n// formatted pyenv(...) example earlier in the document.\nconst target = \"'python.EXE'))\\\"\";\nconst replacement = \"'python.EXE'))\";\n\nconst results = context.document.body.search(target, { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (const range of results.items) {\n  range.insertText(replacement, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The Matlab \"pyenv(...)\" code sample paragraph ended with a stray\n# trailing double-quote character (a leftover from copy/pasting a\n# console transcript): ...'Scripts', 'python.EXE'))\"\n# The revision removes that stray closing quote so the line reads\n# ...'Scripts', 'python.EXE')) -- matching the other, correctly\n# formatted pyenv(...) example earlier in the document.\n\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$searchRange.Find.Text = \"python.EXE'))`\"\"\n$searchRange.Find.Forward = $true\n$searchRange.Find.Wrap = 0\n\n$found = $searchRange.Find.Execute()\n\nif ($found) {\n    # Trim the stray trailing '\"' character directly on the matched\n    # range's text (avoids Find/Replace's smart-quote autocorrect,\n    # which would otherwise mangle the apostrophes in the code sample).\n    $matchedText = $searchRange.Text\n    $searchRange.Text = $matchedText.Substring(0, $matchedText.Length - 1)\n}\n"}
